$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: copy an entire row's values (columns A-R) from $srcRow to $dstRow,
# also copying the NumberFormat of column D (the date column).
# Note: Value2() (rather than Value()) is used so that date-formatted cells
# come back as their raw numeric serial value instead of a DateTime object -
# assigning a DateTime object back into .Value would make Excel silently
# apply a brand new default date NumberFormat to the destination cell.
function Copy-RowValues($srcRow, $dstRow) {
    foreach ($col in @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")) {
        $srcCell = $ws.Range("$col$srcRow")
        $dstCell = $ws.Range("$col$dstRow")
        $dstCell.Value = $srcCell.Value2()
    }
    $ws.Range("D$dstRow").NumberFormat = $ws.Range("D$srcRow").NumberFormat
}

# Helper: copy only the price/date columns (D, J, K, L, M, P) from $srcRow to $dstRow.
function Copy-PriceColumns($srcRow, $dstRow) {
    foreach ($col in @("D","J","K","L","M","P")) {
        $srcCell = $ws.Range("$col$srcRow")
        $dstCell = $ws.Range("$col$dstRow")
        $dstCell.Value = $srcCell.Value2()
    }
}

# 1) Two brand-new rows (446 and 447) are appended at the bottom of the table.
#    They duplicate the current (pre-edit) rows 444 and 445 in full.
Copy-RowValues 444 446
Copy-RowValues 445 447

# 2) Every existing weekly pair of rows from 424..445 takes on the price/date
#    values that used to belong to the pair two rows above it (422..443),
#    effectively pushing all the historical weekly records down by one week
#    (i.e. two rows, since each week occupies a "Primera"/"Segunda" row pair).
#    Processing from the bottom up ensures we always read a row's old values
#    before that row itself gets overwritten.
for ($row = 445; $row -ge 424; $row--) {
    Copy-PriceColumns ($row - 2) $row
}

# 3) Rows 422 and 423 receive the brand new week's data.
$ws.Range("D422").Value = 45041
$ws.Range("J422").Value = 1900

$ws.Range("D423").Value = 45041
$ws.Range("J423").Value = 1000
